{"js": "const replacements = [\n  [\"498\u00f73=166, 0\", \"919\u00f76=153, 1\"],\n  [\"612\u00f73=204, 0\", \"969\u00f72=484, 1\"],\n  [\"140\u00f78=17, 4\", \"855\u00f72=427, 1\"],\n  [\"976\u00f78=122, 0\", \"717\u00f76=119, 3\"],\n  [\"220\u00f76=36, 4\", \"421\u00f73=140, 1\"],\n  [\"197\u00f78=24, 5\", \"397\u00f74=99, 1\"],\n  [\"805\u00f75=161, 0\", \"653\u00f78=81, 5\"],\n  [\"609\u00f77=87, 0\", \"976\u00f74=244, 0\"],\n  [\"818\u00f78=102, 2\", \"894\u00f75=178, 4\"],\n  [\"516\u00f78=64, 4\", \"864\u00f78=108, 0\"],\n  [\"487\u00f73=162, 1\", \"667\u00f72=333, 1\"],\n  [\"877\u00f72=438, 1\", \"171\u00f75=34, 1\"],\n  [\"111\u00f79=12, 3\", \"535\u00f72=267, 1\"],\n  [\"304\u00f74=76, 0\", \"822\u00f79=91, 3\"],\n  [\"965\u00f72=482, 1\", \"760\u00f76=126, 4\"],\n  [\"378\u00f78=47, 2\", \"873\u00f78=109, 1\"],\n  [\"650\u00f74=162, 2\", \"290\u00f74=72, 2\"],\n  [\"800\u00f76=133, 2\", \"398\u00f75=79, 3\"],\n  [\"562\u00f76=93, 4\", \"588\u00f72=294, 0\"],\n  [\"489\u00f79=54, 3\", \"642\u00f75=128, 2\"],\n  [\"504\u00f79=56, 0\", \"717\u00f74=179, 1\"],\n  [\"876\u00f79=97, 3\", \"231\u00f73=77, 0\"],\n  [\"927\u00f79=103, 0\", \"155\u00f73=51, 2\"],\n  [\"564\u00f72=282, 0\", \"899\u00f79=99, 8\"],\n  [\"581\u00f79=64, 5\", \"551\u00f77=78, 5\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('498\u00f73=166, 0', '919\u00f76=153, 1'),\n    @('612\u00f73=204, 0', '969\u00f72=484, 1'),\n    @('140\u00f78=17, 4', '855\u00f72=427, 1'),\n    @('976\u00f78=122, 0', '717\u00f76=119, 3'),\n    @('220\u00f76=36, 4', '421\u00f73=140, 1'),\n    @('197\u00f78=24, 5', '397\u00f74=99, 1'),\n    @('805\u00f75=161, 0', '653\u00f78=81, 5'),\n    @('609\u00f77=87, 0', '976\u00f74=244, 0'),\n    @('818\u00f78=102, 2', '894\u00f75=178, 4'),\n    @('516\u00f78=64, 4', '864\u00f78=108, 0'),\n    @('487\u00f73=162, 1', '667\u00f72=333, 1'),\n    @('877\u00f72=438, 1', '171\u00f75=34, 1'),\n    @('111\u00f79=12, 3', '535\u00f72=267, 1'),\n    @('304\u00f74=76, 0', '822\u00f79=91, 3'),\n    @('965\u00f72=482, 1', '760\u00f76=126, 4'),\n    @('378\u00f78=47, 2', '873\u00f78=109, 1'),\n    @('650\u00f74=162, 2', '290\u00f74=72, 2'),\n    @('800\u00f76=133, 2', '398\u00f75=79, 3'),\n    @('562\u00f76=93, 4', '588\u00f72=294, 0'),\n    @('489\u00f79=54, 3', '642\u00f75=128, 2'),\n    @('504\u00f79=56, 0', '717\u00f74=179, 1'),\n    @('876\u00f79=97, 3', '231\u00f73=77, 0'),\n    @('927\u00f79=103, 0', '155\u00f73=51, 2'),\n    @('564\u00f72=282, 0', '899\u00f79=99, 8'),\n    @('581\u00f79=64, 5', '551\u00f77=78, 5'),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
